# course_data.xlsx — parse/graph prep edit
#
# Changes applied (per the authoritative diff):
#   1. B1: 201 -> 221
#   2. D2:F2 and D3:F3: fill with a single blank-space placeholder value
#      (becomes its own shared string " ").
#   3. The "mon"/"tues" day labels (shared by G1/G2/G3 and J1/J3) are
#      re-cased to "MON"/"TUES" in place.
#   4. The view is re-zoomed to 115% and the active selection moves to L17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- B1: 201 -> 221 ------------------------------------------------------
$ws.Range("B1").Value = 221

# --- Day-label re-casing (mon -> MON, tues -> TUES) -----------------------
# All cells sharing the "mon"/"tues" string are rewritten identically so the
# workbook's shared-string table collapses back onto a single re-cased
# entry per word (same slot reused, no stray duplicate strings left behind).
$ws.Range("G1").Value = "MON"
$ws.Range("G2").Value = "MON"
$ws.Range("G3").Value = "MON"
$ws.Range("J1").Value = "TUES"
$ws.Range("J3").Value = "TUES"

# --- New blank-space cells D2:F2 and D3:F3 --------------------------------
$ws.Range("D2").Value = " "
$ws.Range("E2").Value = " "
$ws.Range("F2").Value = " "
$ws.Range("D3").Value = " "
$ws.Range("E3").Value = " "
$ws.Range("F3").Value = " "

# --- View state: zoom to 115% and move the selection to L17 --------------
$excel.ActiveWindow.Zoom = 115
$ws.Range("L17").Select()
